# Intervjuguide.docx - "minor change in the interviewguide"
#
# 1) The "Tidsramme: ..." paragraph's time estimate changes from
#    "ca en time, kanskje litt mer for noen av intervjuene" to
#    "mellom en og to timer". Word keeps the untouched "Tidsramme: "
#    prefix in its own run and puts the freshly typed replacement text
#    in a new run right after it.
# 2) The "_GoBack" bookmark (Word's "last edit location" marker) moves
#    from the end of the document (right after "Andre?") to the end of
#    this freshly edited text - exactly what real Word does after an
#    edit.

$d = $word.ActiveDocument

# --- Step 1: replace the whole sentence in one go -------------------
# This keeps the run's existing character formatting (<w:lang .../> )
# and correctly marks the trailing-space run as xml:space="preserve".
$whole = $d.Content
$whole.Find.Execute( `
    "Tidsramme: ca en time, kanskje litt mer for noen av intervjuene", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Tidsramme: mellom en og to timer", 2) | Out-Null

# --- Step 2: find the seam between "Tidsramme: " and the new text ---
$prefixRange = $d.Content
$prefixRange.Find.Execute("Tidsramme: ", $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $prefixRange.End

# Briefly drop a bookmark right on that seam: this is enough to make
# the host split the single merged run back into two runs (one for the
# untouched "Tidsramme: " prefix, one for the newly typed text) while
# both keep the original run formatting - matching exactly what Word
# itself leaves behind after a partial-text edit.
$d.Bookmarks.Add("_EditSeamTmp", $d.Range($splitPos, $splitPos)) | Out-Null

# --- Step 3: find the end of the newly typed text --------------------
$newTextRange = $d.Content
$newTextRange.Find.Execute("mellom en og to timer", $true, $false, `
    $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$endPos = $newTextRange.End

$d.Bookmarks.Item("_EditSeamTmp").Delete()

# --- Step 4: move "_GoBack" to the end of the newly typed text -------
# Adding a bookmark named "_GoBack" anywhere relocates the single,
# document-wide "_GoBack" bookmark (removing it from its old spot after
# "Andre?"). Dropping a zero-length bookmark exactly at a paragraph's
# final boundary confuses the host, so a scratch character is
# temporarily appended right after the insertion point, the bookmark is
# planted just before it, and the scratch character is removed again.
$scratch = $d.Range($endPos, $endPos)
$scratch.InsertAfter("X")

$goBackRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

$scratchRange = $d.Range($endPos, $endPos + 1)
$scratchRange.Text = ""
